# feat: add 2022-Q4 data
#
# 1. Insert a new worksheet "2022-Q4" right after "总计", pushing the
#    existing quarter sheets (2022-Q2, 2022-Q1, 2021-Q4, 2021-Q3) one tab
#    to the right (their own data is unchanged).
# 2. Populate "2022-Q4" with the per-fund holdings for that quarter.
# 3. Insert a new row into "总计" (right under the header) with the
#    2022-Q4 totals, shifting the previously-existing rows down by one.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item("总计")
$q2    = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------
# 1) New "2022-Q4" sheet, inserted right after "总计"
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Add($null, $total)
$q4.Name = "2022-Q4"

# Header row - same headers/style as the other per-quarter sheets.
$q4.Cells.Item(1,2).Value = "基金代码"
$q4.Cells.Item(1,3).Value = "基金名称"
$q4.Cells.Item(1,4).Value = "基金规模"
$q4.Cells.Item(1,5).Value = "股票总仓位"
$q4.Cells.Item(1,6).Value = "仓位占比"
$q4.Cells.Item(1,7).Value = "持有市值(亿元)"
$q4.Cells.Item(1,8).Value = "仓位排名"
$q4.Range("B1:H1").Font.Bold = $true
$q4.Range("B1:H1").HorizontalAlignment = -4108
$q4.Range("B1:H1").VerticalAlignment = -4160
$q4.Range("B1:H1").Borders.LineStyle = 1
$q4.Cells.Item(1,1).Font.Bold = $true
$q4.Cells.Item(1,1).HorizontalAlignment = -4108
$q4.Cells.Item(1,1).VerticalAlignment = -4160
$q4.Cells.Item(1,1).Borders.LineStyle = 1

function Set-FundRow {
    param($sheet, $row, $idx, $code, $name, $scale, $pos, $ratio, $mv, $rank)

    $sheet.Cells.Item($row,1).Value = $idx
    $sheet.Cells.Item($row,1).Font.Bold = $true
    $sheet.Cells.Item($row,1).HorizontalAlignment = -4108
    $sheet.Cells.Item($row,1).VerticalAlignment = -4160
    $sheet.Cells.Item($row,1).Borders.LineStyle = 1

    $sheet.Cells.Item($row,2).NumberFormat = "@"
    $sheet.Cells.Item($row,2).Value = $code

    $sheet.Cells.Item($row,3).Value = $name

    $sheet.Cells.Item($row,4).NumberFormat = "@"
    $sheet.Cells.Item($row,4).Value = $scale

    $sheet.Cells.Item($row,5).NumberFormat = "@"
    $sheet.Cells.Item($row,5).Value = $pos

    $sheet.Cells.Item($row,6).NumberFormat = "@"
    $sheet.Cells.Item($row,6).Value = $ratio

    $sheet.Cells.Item($row,7).NumberFormat = "@"
    $sheet.Cells.Item($row,7).Value = $mv

    $sheet.Cells.Item($row,8).Value = $rank
}

Set-FundRow $q4 2 0 "513360" "博时中证全球中国教育主题ETF（QDII）" "4.89" "99.23" "14.63" "0.7154" 2
Set-FundRow $q4 3 1 "000988" "嘉实全球互联网股票-人民币（QDII）" "12.08" "89.83" "5.04" "0.6088" 6
Set-FundRow $q4 4 2 "000989" "嘉实全球互联网股票-美元现汇（QDII）" "12.08" "89.83" "5.04" "0.6088" 6
Set-FundRow $q4 5 3 "000990" "嘉实全球互联网股票-美元现钞（QDII）" "12.08" "89.83" "5.04" "0.6088" 6

# ---------------------------------------------------------------------
# 2) "总计" sheet: insert a new row 2 with the 2022-Q4 totals, pushing
#    the existing rows (2022-Q2 ... 2021-Q3) down by one.
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()

$total.Cells.Item(2,1).Value = 0
$total.Cells.Item(2,1).Font.Bold = $true
$total.Cells.Item(2,1).HorizontalAlignment = -4108
$total.Cells.Item(2,1).VerticalAlignment = -4160
$total.Cells.Item(2,1).Borders.LineStyle = 1

$total.Cells.Item(2,2).Value = "2022-Q4"
$total.Cells.Item(2,3).Value = 4
$total.Cells.Item(2,4).Value = 2.54

# Re-number the index column (A) for the rows that shifted down, and
# restore their original "index" style.
$total.Cells.Item(3,1).Value = 1
$total.Cells.Item(4,1).Value = 2
$total.Cells.Item(5,1).Value = 3
$total.Cells.Item(6,1).Value = 4
"A3","A4","A5","A6" | ForEach-Object {
    $c = $total.Range($_)
    $c.Font.Bold = $true
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
    $c.Borders.LineStyle = 1
}

# Restore the original active tab ("总计" is tab 0, same as before the edit).
$total.Select()
$total.Range("A1").Select()
